# Insert two missing historical daily-price rows (2019-11-18 and 2019-11-26)
# into the "0188" (HLT) price history sheet, right before the existing
# 2019-11-29 row. Everything below shifts down by two rows, which Excel's
# native row-insert handles for us (including auto-growing the used range /
# <dimension>).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 696/697 (existing data - previously starting
# with 2019-11-29 at row 696 - shifts down to rows 698+).
$ws.Rows.Item(696).Insert()
$ws.Rows.Item(697).Insert()

# Helper-free, explicit cell-by-cell population so every column keeps the
# exact same literal type the rest of the sheet uses:
#   A = unix timestamp (number), B = date (text), C = id (text, "0188"),
#   D = name (text, "HLT"), E-H = price (number), I = volume (number, or the
#   literal text "-" when there was no trading that day).

function Set-TextCell($cell, [string]$text) {
    # Force literal text so numeric-looking strings ("0188", "2019-11-18")
    # are not reinterpreted as a number/date, then drop the temporary
    # Text number-format so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 696: 2019-11-18
$r = $ws.Cells.Item(696, 1); $r.Value = 1574035200
Set-TextCell $ws.Cells.Item(696, 2) "2019-11-18"
Set-TextCell $ws.Cells.Item(696, 3) "0188"
$ws.Cells.Item(696, 4).Value = "HLT"
$ws.Cells.Item(696, 5).Value = 0.185
$ws.Cells.Item(696, 6).Value = 0.19
$ws.Cells.Item(696, 7).Value = 0.185
$ws.Cells.Item(696, 8).Value = 0.19
$ws.Cells.Item(696, 9).Value = 171000

# Row 697: 2019-11-26
$ws.Cells.Item(697, 1).Value = 1574726400
Set-TextCell $ws.Cells.Item(697, 2) "2019-11-26"
Set-TextCell $ws.Cells.Item(697, 3) "0188"
$ws.Cells.Item(697, 4).Value = "HLT"
$ws.Cells.Item(697, 5).Value = 0.19
$ws.Cells.Item(697, 6).Value = 0.19
$ws.Cells.Item(697, 7).Value = 0.19
$ws.Cells.Item(697, 8).Value = 0.19
$ws.Cells.Item(697, 9).Value = 100000
